# Auto-generated edit script: updates market-price/profit calculation cells
# (columns H-N) across several rows in the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets,
# matching a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 500
$ws.Range("I33").Value = 269.14285
$ws.Range("J33").Value = 1038.6666
$ws.Range("K33").Value = 269.14285
$ws.Range("L33").Value = 1038.6666
$ws.Range("M33").Value = -40.14285000000001
$ws.Range("N33").Value = -1496.6666

$ws.Range("H41").Value = 84.333336
$ws.Range("I41").Value = 87
$ws.Range("J41").Value = 79
$ws.Range("K41").Value = 87
$ws.Range("L41").Value = 79
$ws.Range("M41").Value = 353
$ws.Range("N41").Value = -959

$ws.Range("H55").Value = 331.91666
$ws.Range("J55").Value = 459.5
$ws.Range("L55").Value = 459.5
$ws.Range("N55").Value = -887.5

$ws.Range("H80").Value = 8334593.5
$ws.Range("I80").Value = 13889700
$ws.Range("J80").Value = 1933.5
$ws.Range("K80").Value = 41669100
$ws.Range("L80").Value = 5800.5
$ws.Range("M80").Value = -41668102
$ws.Range("N80").Value = -7796.5

$ws.Range("H83").Value = 8334593.5
$ws.Range("I83").Value = 13889700
$ws.Range("J83").Value = 1933.5
$ws.Range("K83").Value = 125007300
$ws.Range("L83").Value = 17401.5
$ws.Range("M83").Value = -125002308
$ws.Range("N83").Value = -27385.5

$ws.Range("H132").Value = 82271.516
$ws.Range("I132").Value = 82271.516
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 246814.548
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -244284.548
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 1302002.6
$ws.Range("I137").Value = 3444.9092
$ws.Range("K137").Value = 10334.7276
$ws.Range("M137").Value = -7784.7276

$ws.Range("H138").Value = 1710.1282
$ws.Range("I138").Value = 1042.4073
$ws.Range("J138").Value = 3212.5
$ws.Range("K138").Value = 3127.2219
$ws.Range("L138").Value = 9637.5
$ws.Range("M138").Value = 2012.7781
$ws.Range("N138").Value = -19917.5

$ws.Range("H141").Value = 1763
$ws.Range("I141").Value = 1580.1818
$ws.Range("K141").Value = 4740.5454
$ws.Range("M141").Value = 439.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11431
$ws.Range("J2").Value = 17135.111
$ws.Range("L2").Value = 17135.111
$ws.Range("N2").Value = -17361.111

$ws.Range("H28").Value = 1731.4286
$ws.Range("I28").Value = 1731.4286
$ws.Range("K28").Value = 1731.4286
$ws.Range("M28").Value = -1539.4286

$ws.Range("H32").Value = 6292423
$ws.Range("I32").Value = 6669944.5
$ws.Range("K32").Value = 6669944.5
$ws.Range("M32").Value = -6669657.5

$ws.Range("H45").Value = 2473.125
$ws.Range("I45").Value = 2512.1428
$ws.Range("J45").Value = 2200
$ws.Range("K45").Value = 2512.1428
$ws.Range("L45").Value = 2200
$ws.Range("M45").Value = -2135.1428
$ws.Range("N45").Value = -2954

$ws.Range("H88").Value = 1206.6154
$ws.Range("I88").Value = 762.5
$ws.Range("J88").Value = 1404
$ws.Range("K88").Value = 762.5
$ws.Range("L88").Value = 1404
$ws.Range("M88").Value = -356.5
$ws.Range("N88").Value = -2216

$ws.Range("H91").Value = 1206.6154
$ws.Range("I91").Value = 762.5
$ws.Range("J91").Value = 1404
$ws.Range("K91").Value = 762.5
$ws.Range("L91").Value = 1404
$ws.Range("M91").Value = 641.5
$ws.Range("N91").Value = -4212

$ws.Range("H99").Value = 1731.4286
$ws.Range("I99").Value = 1731.4286
$ws.Range("K99").Value = 1731.4286
$ws.Range("M99").Value = 1263.5714

$ws.Range("H116").Value = 11431
$ws.Range("J116").Value = 17135.111
$ws.Range("L116").Value = 17135.111
$ws.Range("N116").Value = -21723.111

$ws.Range("H132").Value = 1696505.2
$ws.Range("I132").Value = 1978589.4
$ws.Range("K132").Value = 5935768.199999999
$ws.Range("M132").Value = -5933238.199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11431
$ws.Range("J3").Value = 17135.111
$ws.Range("L3").Value = 17135.111
$ws.Range("N3").Value = -17363.111

$ws.Range("H86").Value = 962
$ws.Range("I86").Value = 840.6667
$ws.Range("J86").Value = 1144
$ws.Range("K86").Value = 840.6667
$ws.Range("L86").Value = 1144
$ws.Range("M86").Value = 282.3333
$ws.Range("N86").Value = -3390

$ws.Range("H89").Value = 962
$ws.Range("I89").Value = 840.6667
$ws.Range("J89").Value = 1144
$ws.Range("K89").Value = 4203.3335
$ws.Range("L89").Value = 5720
$ws.Range("M89").Value = 1412.6665
$ws.Range("N89").Value = -16952

$ws.Range("H134").Value = 2030675.4
$ws.Range("I134").Value = 2384603.2
$ws.Range("K134").Value = 7153809.600000001
$ws.Range("M134").Value = -7151274.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 117368.1
$ws.Range("I31").Value = 157280.9
$ws.Range("J31").Value = 37542.5
$ws.Range("K31").Value = 157280.9
$ws.Range("L31").Value = 37542.5
$ws.Range("M31").Value = -156985.9
$ws.Range("N31").Value = -38132.5

$ws.Range("H34").Value = 117368.1
$ws.Range("I34").Value = 157280.9
$ws.Range("J34").Value = 37542.5
$ws.Range("K34").Value = 157280.9
$ws.Range("L34").Value = 37542.5
$ws.Range("M34").Value = -157078.9
$ws.Range("N34").Value = -37946.5

$ws.Range("H62").Value = 3580
$ws.Range("J62").Value = 4800
$ws.Range("L62").Value = 4800
$ws.Range("N62").Value = -6048

$ws.Range("H65").Value = 3580
$ws.Range("J65").Value = 4800
$ws.Range("L65").Value = 24000
$ws.Range("N65").Value = -30240

$ws.Range("H132").Value = 225036.56
$ws.Range("I132").Value = 2135.9
$ws.Range("J132").Value = 1711041
$ws.Range("K132").Value = 6407.700000000001
$ws.Range("L132").Value = 5133123
$ws.Range("M132").Value = -3877.700000000001
$ws.Range("N132").Value = -5138183

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 272.375
$ws.Range("I8").Value = 272.375
$ws.Range("K8").Value = 817.125
$ws.Range("M8").Value = -678.125

$ws.Range("H23").Value = 114.53333
$ws.Range("I23").Value = 52.545456
$ws.Range("K23").Value = 157.636368
$ws.Range("M23").Value = 77.363632

$ws.Range("H123").Value = 2999.8572
$ws.Range("I123").Value = 2999.8572
$ws.Range("K123").Value = 8999.571599999999
$ws.Range("M123").Value = -6549.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 32499.8
$ws.Range("I20").Value = 29500
$ws.Range("J20").Value = 33249.75
$ws.Range("K20").Value = 29500
$ws.Range("L20").Value = 33249.75
$ws.Range("M20").Value = -29255
$ws.Range("N20").Value = -33739.75

$ws.Range("H132").Value = 812773.6
$ws.Range("I132").Value = 1213902
$ws.Range("J132").Value = 10516.8
$ws.Range("K132").Value = 3641706
$ws.Range("L132").Value = 31550.4
$ws.Range("M132").Value = -3639176
$ws.Range("N132").Value = -36610.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1721.5555
$ws.Range("I22").Value = 998.8
$ws.Range("K22").Value = 998.8
$ws.Range("M22").Value = -703.8

$ws.Range("H27").Value = 1721.5555
$ws.Range("I27").Value = 998.8
$ws.Range("K27").Value = 998.8
$ws.Range("M27").Value = -891.8

$ws.Range("H82").Value = 197
$ws.Range("I82").Value = 196
$ws.Range("J82").Value = 199
$ws.Range("K82").Value = 196
$ws.Range("L82").Value = 199
$ws.Range("M82").Value = 165
$ws.Range("N82").Value = -921

$ws.Range("H85").Value = 197
$ws.Range("I85").Value = 196
$ws.Range("J85").Value = 199
$ws.Range("K85").Value = 196
$ws.Range("L85").Value = 199
$ws.Range("M85").Value = 1052
$ws.Range("N85").Value = -2695

$ws.Range("H132").Value = 2681087
$ws.Range("I132").Value = 4973306.5
$ws.Range("K132").Value = 14919919.5
$ws.Range("M132").Value = -14917389.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1991.3334
$ws.Range("I81").Value = 1991.3334
$ws.Range("K81").Value = 3982.6668
$ws.Range("M81").Value = -2921.6668

$ws.Range("H84").Value = 1991.3334
$ws.Range("I84").Value = 1991.3334
$ws.Range("K84").Value = 19913.334
$ws.Range("M84").Value = -14609.334
